$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad" / Changed date) for rows 2-15 from 45182 (2023-09-13)
# to 45184 (2023-09-15), keeping existing number formatting/style intact.
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
